# Add a new "Abstract Title" paragraph style and tighten the spacing
# above the existing "Abstract" style, matching the authored diff:
#   - new custom style "AbstractTitle" ("Abstract Title"), based on
#     Normal, followed by Abstract, centered, bold, small, blue text,
#     with keepNext/keepLines and spacing before=300(twips)/after=0.
#   - Abstract style: spacing-before reduced from 300 -> 100 (twentieths
#     of a point), i.e. from 15pt to 5pt. SpaceAfter (300 = 15pt) is
#     left unchanged.

$d = $word.ActiveDocument

# --- 1. Create the "Abstract Title" style -----------------------------
$abstractTitle = $d.Styles.Add("Abstract Title", 1)  # 1 = wdStyleTypeParagraph

$abstractTitle.BaseStyle = $d.Styles("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles("Abstract")
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1          # wdAlignParagraphCenter
$abstractTitle.ParagraphFormat.SpaceBefore = 15        # 300 twentieths-of-a-point
$abstractTitle.ParagraphFormat.SpaceAfter = 0

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060                    # RGB(0x34,0x5A,0x8A) -> 345A8A

# --- 2. Tighten the spacing above the "Abstract" style -----------------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5               # 100 twentieths-of-a-point

Write-Output "AbstractTitle style added; Abstract spacing-before updated."
